# Generate Report for Handoff
# Updates the "b.md" row across the Overview / zh-cn / de-de sheets to reflect
# that a new handoff package (b.63290e5768f688058c7b37413b0a5c26c308f864)
# was generated and is now "Ready for handoff".

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the "b.md" file ---
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("B3").Value = "Ready for handoff"
$ovw.Range("C3").Value = "Ready for handoff"
$ovw.Range("D3").Value = "2016-03-18 16:48:11"

# --- zh-cn sheet: row 3 is the "b.md" file ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("E3").Value = "2016-03-18 16:48:01"

# --- de-de sheet: row 3 is the "b.md" file ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("E3").Value = "2016-03-18 16:48:11"
